# Finalize the trade recorded in row 2: it closed out, so fill in the
# previously-blank sell details and mark the trade as no longer held /
# not profitable.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = $false                  # Profitable
$ws.Range("E2").Value = 309.77999999999997       # SellPrice
$ws.Range("F2").Value = -0.6159769008662227      # Price Change %
$ws.Range("G2").Value = $false                   # Holding (no longer held)

# New trade row logged after closing the previous one.
$ws.Range("C3").Value = 9938.4                   # Principle
